$wb = $excel.ActiveWorkbook

# The trading bot logged a new open MarketMaking trade (#30) appended as
# row 31 on both the "All Trades" and "MarketMaking" worksheets.
$sheetNames = @("All Trades", "MarketMaking")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $row = 31

    # Plain numeric columns
    $ws.Range("A$row").Value = 30
    $ws.Range("F$row").Value = 0.96
    $ws.Range("I$row").Value = 0
    $ws.Range("J$row").Value = 0
    $ws.Range("K$row").Value = 99.24617257389608
    $ws.Range("L$row").Value = 0
    $ws.Range("M$row").Value = 0
    $ws.Range("N$row").Value = 0.6
    $ws.Range("Q$row").Value = 0

    # Date column: force text so Excel doesn't coerce "2026-02-17" into a
    # date serial number, then drop the temporary number format again so
    # no stray style sticks to the cell.
    $ws.Range("B$row").NumberFormat = "@"
    $ws.Range("B$row").Value = "2026-02-17"
    $ws.Range("B$row").Style = "Normal"

    # Plain text columns
    $ws.Range("C$row").Value = "13:19:24"
    $ws.Range("D$row").Value = "MarketMaking"
    $ws.Range("E$row").Value = "DOWN"
    $ws.Range("H$row").Value = "OPEN"
    $ws.Range("O$row").Value = "Normal spread capture: 19600 bps"

    # Exit price / exit reason are blank because the trade is still OPEN.
    # A leading quote forces an explicit (empty) text cell instead of
    # Excel silently dropping the write for a plain "" assignment; reset
    # the style afterwards so no quote-prefix formatting lingers.
    $ws.Range("G$row").Value = "'"
    $ws.Range("G$row").Style = "Normal"
    $ws.Range("P$row").Value = "'"
    $ws.Range("P$row").Style = "Normal"
}
